# Automatic update: refresh the "Förändrad" (Last changed) date stamp
# held in column C for every data row of the sheet. The tracked date
# advances from serial 45172 (2023-09-03) to serial 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column C (falls back to the
# sheet's UsedRange if column C turns out to be empty).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

$newValue = 45175

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = $newValue
    }
}
